$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 752.7917
$ws.Range("J17").Value = 752.7917
$ws.Range("L17").Value = 2258.3751
$ws.Range("N17").Value = -2594.3751

# Row 58
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 9000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -9300

# Row 92
$ws.Range("H92").Value = 1551.4
$ws.Range("I92").Value = 622.8333
$ws.Range("J92").Value = 3939.1428
$ws.Range("K92").Value = 622.8333
$ws.Range("L92").Value = 3939.1428
$ws.Range("M92").Value = 625.1667
$ws.Range("N92").Value = -6435.1428

# Row 94
$ws.Range("H94").Value = 998.6667
$ws.Range("I94").Value = 998.6667
$ws.Range("K94").Value = 998.6667
$ws.Range("M94").Value = -547.6667

# Row 96
$ws.Range("H96").Value = 699.3
$ws.Range("I96").Value = 680.2
$ws.Range("J96").Value = 718.4
$ws.Range("K96").Value = 2040.6
$ws.Range("L96").Value = 2155.2
$ws.Range("M96").Value = -667.6000000000001
$ws.Range("N96").Value = -4901.2

# Row 100
$ws.Range("H100").Value = 4517.6665
$ws.Range("I100").Value = 1072.8889
$ws.Range("J100").Value = 7962.4443
$ws.Range("K100").Value = 1072.8889
$ws.Range("L100").Value = 7962.4443
$ws.Range("M100").Value = -531.8888999999999
$ws.Range("N100").Value = -9044.444299999999

# Row 107
$ws.Range("H107").Value = 1555.7142
$ws.Range("I107").Value = 1548
$ws.Range("K107").Value = 1548
$ws.Range("M107").Value = 372

# Row 113
$ws.Range("H113").Value = 6732.482
$ws.Range("I113").Value = 6643.4326
$ws.Range("K113").Value = 6643.4326
$ws.Range("M113").Value = -3389.4326

# Row 138
$ws.Range("H138").Value = 3866.9412
$ws.Range("I138").Value = 4668.8335
$ws.Range("J138").Value = 3429.5454
$ws.Range("K138").Value = 14006.5005
$ws.Range("L138").Value = 10288.6362
$ws.Range("M138").Value = -8866.500499999998
$ws.Range("N138").Value = -20568.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 119
$ws.Range("H119").Value = 19998
$ws.Range("J119").Value = 19998
$ws.Range("L119").Value = 19998
$ws.Range("N119").Value = -29674

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1099.4
$ws.Range("J7").Value = 999
$ws.Range("L7").Value = 999
$ws.Range("N7").Value = -1225

# Row 20
$ws.Range("H20").Value = 2240
$ws.Range("I20").Value = 1500
$ws.Range("J20").Value = 2733.3333
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 2733.3333
$ws.Range("M20").Value = -1253
$ws.Range("N20").Value = -3227.3333

# Row 105
$ws.Range("H105").Value = 1647.8125
$ws.Range("I105").Value = 1690.2667
$ws.Range("K105").Value = 1690.2667
$ws.Range("M105").Value = 56.7333000000001

# Row 134
$ws.Range("H134").Value = 7104
$ws.Range("I134").Value = 5853.8066
$ws.Range("J134").Value = 12640.571
$ws.Range("K134").Value = 17561.4198
$ws.Range("L134").Value = 37921.713
$ws.Range("M134").Value = -15026.4198
$ws.Range("N134").Value = -42991.713

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 34749
$ws.Range("I51").Value = 14499.5
$ws.Range("J51").Value = 54998.5
$ws.Range("K51").Value = 14499.5
$ws.Range("L51").Value = 54998.5
$ws.Range("M51").Value = -13763.5
$ws.Range("N51").Value = -56470.5

# Row 61
$ws.Range("H61").Value = 34749
$ws.Range("I61").Value = 14499.5
$ws.Range("J61").Value = 54998.5
$ws.Range("K61").Value = 14499.5
$ws.Range("L61").Value = 54998.5
$ws.Range("M61").Value = -14151.5
$ws.Range("N61").Value = -55694.5

# Row 62
$ws.Range("H62").Value = 10946.5
$ws.Range("I62").Value = 4031
$ws.Range("J62").Value = 28235.25
$ws.Range("K62").Value = 4031
$ws.Range("L62").Value = 28235.25
$ws.Range("M62").Value = -3407
$ws.Range("N62").Value = -29483.25

# Row 65
$ws.Range("H65").Value = 10946.5
$ws.Range("I65").Value = 4031
$ws.Range("J65").Value = 28235.25
$ws.Range("K65").Value = 20155
$ws.Range("L65").Value = 141176.25
$ws.Range("M65").Value = -17035
$ws.Range("N65").Value = -147416.25

# Row 105
$ws.Range("H105").Value = 1347.8334
$ws.Range("I105").Value = 1247.7142
$ws.Range("J105").Value = 1698.25
$ws.Range("K105").Value = 1247.7142
$ws.Range("L105").Value = 1698.25
$ws.Range("M105").Value = 499.2858000000001
$ws.Range("N105").Value = -5192.25

# Row 132
$ws.Range("H132").Value = 5931.6665
$ws.Range("J132").Value = 3644.5
$ws.Range("L132").Value = 10933.5
$ws.Range("N132").Value = -15993.5

# Row 134
$ws.Range("H134").Value = 10844.777
$ws.Range("J134").Value = 14336.333
$ws.Range("L134").Value = 43008.999
$ws.Range("N134").Value = -48078.999

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 461.5
$ws.Range("J23").Value = 490.75
$ws.Range("L23").Value = 1472.25
$ws.Range("N23").Value = -1942.25

# Row 58
$ws.Range("H58").Value = 9500.5
$ws.Range("J58").Value = 10000.6
$ws.Range("L58").Value = 30001.8
$ws.Range("N58").Value = -30257.8

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 519.2
$ws.Range("I2").Value = 499
$ws.Range("K2").Value = 499
$ws.Range("M2").Value = -386

# Row 4
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724

# Row 70
$ws.Range("H70").Value = 7327.25
$ws.Range("I70").Value = 6885.1177
$ws.Range("K70").Value = 6885.1177
$ws.Range("M70").Value = -6615.1177

# Row 73
$ws.Range("H73").Value = 7327.25
$ws.Range("I73").Value = 6885.1177
$ws.Range("K73").Value = 6885.1177
$ws.Range("M73").Value = -5949.1177

# Row 97
$ws.Range("H97").Value = 835.1667
$ws.Range("I97").Value = 328.93332
$ws.Range("K97").Value = 328.93332
$ws.Range("M97").Value = 167.06668

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 4499.5
$ws.Range("I93").Value = 2998
$ws.Range("K93").Value = 2998
$ws.Range("M93").Value = -1750

# Row 136
$ws.Range("I136").Value = 1250.7059
$ws.Range("J136").Value = 3586.5454
$ws.Range("K136").Value = 3752.1177
$ws.Range("L136").Value = 10759.6362
$ws.Range("M136").Value = -1202.1177
$ws.Range("N136").Value = -15859.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 43865.75
$ws.Range("J52").Value = 95000
$ws.Range("L52").Value = 95000
$ws.Range("N52").Value = -95452

# Row 100
$ws.Range("H100").Value = 400.14285
$ws.Range("I100").Value = 333.45456
$ws.Range("K100").Value = 666.90912
$ws.Range("M100").Value = -125.90912

# Row 107
$ws.Range("H107").Value = 5913.9473
$ws.Range("I107").Value = 6201.4443
$ws.Range("K107").Value = 18604.3329
$ws.Range("M107").Value = -16684.3329
